# Add a "beginner Examination range" list of LeetCode problem numbers to
# column A, starting at row 7 (rows 1-6 already contain the header/table),
# continuing down to row 67.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @(
    3,4,5,11,14,16,19,20,21,24,
    27,28,34,35,54,56,58,66,67,74,
    83,92,136,151,153,155,167,169,189,190,
    191,201,205,209,217,231,234,242,258,268,
    338,342,349,405,461,476,485,496,500,557,
    561,621,622,709,724,905,922,933,961,1013,
    1207
)

$startRow = 7
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $values[$i]
}

$lastRow = $startRow + $values.Length - 1

# Match the saved view state from the diff: zoomed to 125%, with A67 selected.
$ws.Range("A" + $lastRow).Select()
$excel.ActiveWindow.Zoom = 125
